{"js": "const replacements = [\n  [\"2025-02-24 Monday\", \"2025-02-25 Tuesday\"],\n  [\"604\u00f77=86, 2\", \"686\u00f79=76, 2\"],\n  [\"357\u00f78=44, 5\", \"576\u00f74=144, 0\"],\n  [\"407\u00f73=135, 2\", \"253\u00f77=36, 1\"],\n  [\"149\u00f75=29, 4\", \"497\u00f72=248, 1\"],\n  [\"231\u00f78=28, 7\", \"470\u00f78=58, 6\"],\n  [\"624\u00f78=78, 0\", \"826\u00f76=137, 4\"],\n  [\"538\u00f75=107, 3\", \"161\u00f73=53, 2\"],\n  [\"675\u00f73=225, 0\", \"131\u00f72=65, 1\"],\n  [\"810\u00f78=101, 2\", \"788\u00f79=87, 5\"],\n  [\"714\u00f78=89, 2\", \"346\u00f77=49, 3\"],\n  [\"721\u00f78=90, 1\", \"598\u00f73=199, 1\"],\n  [\"838\u00f73=279, 1\", \"818\u00f74=204, 2\"],\n  [\"884\u00f78=110, 4\", \"224\u00f78=28, 0\"],\n  [\"906\u00f77=129, 3\", \"962\u00f76=160, 2\"],\n  [\"746\u00f74=186, 2\", \"480\u00f75=96, 0\"],\n  [\"390\u00f79=43, 3\", \"239\u00f79=26, 5\"],\n  [\"467\u00f77=66, 5\", \"836\u00f79=92, 8\"],\n  [\"407\u00f72=203, 1\", \"731\u00f73=243, 2\"],\n  [\"112\u00f75=22, 2\", \"782\u00f73=260, 2\"],\n  [\"130\u00f77=18, 4\", \"446\u00f73=148, 2\"],\n  [\"929\u00f74=232, 1\", \"779\u00f73=259, 2\"],\n  [\"949\u00f79=105, 4\", \"968\u00f78=121, 0\"],\n  [\"945\u00f73=315, 0\", \"207\u00f78=25, 7\"],\n  [\"988\u00f77=141, 1\", \"680\u00f73=226, 2\"],\n  [\"235\u00f78=29, 3\", \"629\u00f75=125, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = '2025-02-24 Monday'; New = '2025-02-25 Tuesday' },\n  @{ Old = '604\u00f77=86, 2'; New = '686\u00f79=76, 2' },\n  @{ Old = '357\u00f78=44, 5'; New = '576\u00f74=144, 0' },\n  @{ Old = '407\u00f73=135, 2'; New = '253\u00f77=36, 1' },\n  @{ Old = '149\u00f75=29, 4'; New = '497\u00f72=248, 1' },\n  @{ Old = '231\u00f78=28, 7'; New = '470\u00f78=58, 6' },\n  @{ Old = '624\u00f78=78, 0'; New = '826\u00f76=137, 4' },\n  @{ Old = '538\u00f75=107, 3'; New = '161\u00f73=53, 2' },\n  @{ Old = '675\u00f73=225, 0'; New = '131\u00f72=65, 1' },\n  @{ Old = '810\u00f78=101, 2'; New = '788\u00f79=87, 5' },\n  @{ Old = '714\u00f78=89, 2'; New = '346\u00f77=49, 3' },\n  @{ Old = '721\u00f78=90, 1'; New = '598\u00f73=199, 1' },\n  @{ Old = '838\u00f73=279, 1'; New = '818\u00f74=204, 2' },\n  @{ Old = '884\u00f78=110, 4'; New = '224\u00f78=28, 0' },\n  @{ Old = '906\u00f77=129, 3'; New = '962\u00f76=160, 2' },\n  @{ Old = '746\u00f74=186, 2'; New = '480\u00f75=96, 0' },\n  @{ Old = '390\u00f79=43, 3'; New = '239\u00f79=26, 5' },\n  @{ Old = '467\u00f77=66, 5'; New = '836\u00f79=92, 8' },\n  @{ Old = '407\u00f72=203, 1'; New = '731\u00f73=243, 2' },\n  @{ Old = '112\u00f75=22, 2'; New = '782\u00f73=260, 2' },\n  @{ Old = '130\u00f77=18, 4'; New = '446\u00f73=148, 2' },\n  @{ Old = '929\u00f74=232, 1'; New = '779\u00f73=259, 2' },\n  @{ Old = '949\u00f79=105, 4'; New = '968\u00f78=121, 0' },\n  @{ Old = '945\u00f73=315, 0'; New = '207\u00f78=25, 7' },\n  @{ Old = '988\u00f77=141, 1'; New = '680\u00f73=226, 2' },\n  @{ Old = '235\u00f78=29, 3'; New = '629\u00f75=125, 4' },\n)\n\nforeach ($r in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $ok = $find.Execute(\n    [ref]$r.Old,\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]$r.New,\n    [ref]2\n  )\n  if (-not $ok) {\n    Write-Output (\"NOT FOUND: \" + $r.Old)\n  }\n}\n\nWrite-Output \"done\"\n"}
